$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "MuSCs"
$ws.Cells.Item(2,2).Value = "Ifng"
$ws.Cells.Item(2,3).Value = "Ifngr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.06628299999999999
$ws.Cells.Item(2,8).Value = 0.198849
$ws.Cells.Item(2,9).Value = 0.581839834503261
$ws.Cells.Item(2,10).Value = 0.581839834503261
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 7.486252333333333
$ws.Cells.Item(2,14).Value = 22.458757
$ws.Cells.Item(2,15).Value = 0.1371548009319641
$ws.Cells.Item(2,16).Value = 0.1371548009319641
$ws.Cells.Item(2,17).Value = 0.4962112634103333
$ws.Cells.Item(2,18).Value = 4.465901370693
$ws.Cells.Item(2,19).Value = 0.07980212667558169
$ws.Cells.Item(2,20).Value = 0.07980212667558172

# Row 3
$ws.Cells.Item(3,1).Value = "MuSCs"
$ws.Cells.Item(3,2).Value = "Ifng"
$ws.Cells.Item(3,3).Value = "Ifngr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.06628299999999999
$ws.Cells.Item(3,8).Value = 0.198849
$ws.Cells.Item(3,9).Value = 0.581839834503261
$ws.Cells.Item(3,10).Value = 0.581839834503261
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 12.08386866666667
$ws.Cells.Item(3,14).Value = 36.251606
$ws.Cells.Item(3,15).Value = 0.2213872212248432
$ws.Cells.Item(3,16).Value = 0.2213872212248432
$ws.Cells.Item(3,17).Value = 0.8009550668326665
$ws.Cells.Item(3,18).Value = 7.208595601493999
$ws.Cells.Item(3,19).Value = 0.1288119041585996
$ws.Cells.Item(3,20).Value = 0.1288119041585996

# Row 4
$ws.Cells.Item(4,1).Value = "MuSCs"
$ws.Cells.Item(4,2).Value = "Ifng"
$ws.Cells.Item(4,3).Value = "Ifngr2"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.06628299999999999
$ws.Cells.Item(4,8).Value = 0.198849
$ws.Cells.Item(4,9).Value = 0.581839834503261
$ws.Cells.Item(4,10).Value = 0.581839834503261
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 7.866546666666667
$ws.Cells.Item(4,14).Value = 23.59964
$ws.Cells.Item(4,15).Value = 0.1441221313479645
$ws.Cells.Item(4,16).Value = 0.1441221313479645
$ws.Cells.Item(4,17).Value = 0.5214183127066666
$ws.Cells.Item(4,18).Value = 4.69276481436
$ws.Cells.Item(4,19).Value = 0.08385599705175692
$ws.Cells.Item(4,20).Value = 0.08385599705175693

# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Ifng"
$ws.Cells.Item(5,3).Value = "Ifngr2"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.06628299999999999
$ws.Cells.Item(5,8).Value = 0.198849
$ws.Cells.Item(5,9).Value = 0.581839834503261
$ws.Cells.Item(5,10).Value = 0.581839834503261
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 27.14583533333333
$ws.Cells.Item(5,14).Value = 81.437506
$ws.Cells.Item(5,15).Value = 0.4973358464952282
$ws.Cells.Item(5,16).Value = 0.4973358464952283
$ws.Cells.Item(5,17).Value = 1.799307403399333
$ws.Cells.Item(5,18).Value = 16.193766630594
$ws.Cells.Item(5,19).Value = 0.2893698066173228
$ws.Cells.Item(5,20).Value = 0.2893698066173229

# Row 6
$ws.Cells.Item(6,1).Value = "Resolving-Mac"
$ws.Cells.Item(6,2).Value = "Ifng"
$ws.Cells.Item(6,3).Value = "Ifngr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.04763666666666667
$ws.Cells.Item(6,8).Value = 0.14291
$ws.Cells.Item(6,9).Value = 0.4181601654967389
$ws.Cells.Item(6,10).Value = 0.4181601654967389
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 7.486252333333333
$ws.Cells.Item(6,14).Value = 22.458757
$ws.Cells.Item(6,15).Value = 0.1371548009319641
$ws.Cells.Item(6,16).Value = 0.1371548009319641
$ws.Cells.Item(6,17).Value = 0.3566201069855555
$ws.Cells.Item(6,18).Value = 3.20958096287
$ws.Cells.Item(6,19).Value = 0.05735267425638238
$ws.Cells.Item(6,20).Value = 0.05735267425638239

# Row 7
$ws.Cells.Item(7,1).Value = "Resolving-Mac"
$ws.Cells.Item(7,2).Value = "Ifng"
$ws.Cells.Item(7,3).Value = "Ifngr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.04763666666666667
$ws.Cells.Item(7,8).Value = 0.14291
$ws.Cells.Item(7,9).Value = 0.4181601654967389
$ws.Cells.Item(7,10).Value = 0.4181601654967389
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 12.08386866666667
$ws.Cells.Item(7,14).Value = 36.251606
$ws.Cells.Item(7,15).Value = 0.2213872212248432
$ws.Cells.Item(7,16).Value = 0.2213872212248432
$ws.Cells.Item(7,17).Value = 0.5756352237177778
$ws.Cells.Item(7,18).Value = 5.18071701346
$ws.Cells.Item(7,19).Value = 0.09257531706624357
$ws.Cells.Item(7,20).Value = 0.09257531706624357

# Row 8
$ws.Cells.Item(8,1).Value = "Resolving-Mac"
$ws.Cells.Item(8,2).Value = "Ifng"
$ws.Cells.Item(8,3).Value = "Ifngr2"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.04763666666666667
$ws.Cells.Item(8,8).Value = 0.14291
$ws.Cells.Item(8,9).Value = 0.4181601654967389
$ws.Cells.Item(8,10).Value = 0.4181601654967389
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 7.866546666666667
$ws.Cells.Item(8,14).Value = 23.59964
$ws.Cells.Item(8,15).Value = 0.1441221313479645
$ws.Cells.Item(8,16).Value = 0.1441221313479645
$ws.Cells.Item(8,17).Value = 0.3747360613777778
$ws.Cells.Item(8,18).Value = 3.3726245524
$ws.Cells.Item(8,19).Value = 0.06026613429620758
$ws.Cells.Item(8,20).Value = 0.06026613429620759

# Row 9
$ws.Cells.Item(9,1).Value = "Resolving-Mac"
$ws.Cells.Item(9,2).Value = "Ifng"
$ws.Cells.Item(9,3).Value = "Ifngr2"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.04763666666666667
$ws.Cells.Item(9,8).Value = 0.14291
$ws.Cells.Item(9,9).Value = 0.4181601654967389
$ws.Cells.Item(9,10).Value = 0.4181601654967389
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 27.14583533333333
$ws.Cells.Item(9,14).Value = 81.437506
$ws.Cells.Item(9,15).Value = 0.4973358464952282
$ws.Cells.Item(9,16).Value = 0.4973358464952283
$ws.Cells.Item(9,17).Value = 1.293137109162222
$ws.Cells.Item(9,18).Value = 11.63823398246
$ws.Cells.Item(9,19).Value = 0.2079660398779054
$ws.Cells.Item(9,20).Value = 0.2079660398779054
